$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.205933928489685
$ws.Range("B1").Value = 2.326450347900391
$ws.Range("C1").Value = 3.344620943069458
$ws.Range("D1").Value = 3.438175916671753
$ws.Range("E1").Value = 1.133975863456726
